# Apply TSCA2019 culture folders review updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TSCA cultures")

# Row 7 (culture G6): expand on the 72-8208/72-8209 offset note, attribute the offset
# explanation to Plate 72-8208, and revise the "rest time" follow-up note wording.
$ws.Range("C7").Value = 'Regarding 40 min Treated Start time: Plate 72-8208 "hit stop began recording while offset was running"
Regarding 40 min Treated Start time: Plate 72-8209 "recording had stopped. Thought I had forgot to start it. Started it again then realized mistake. Actual start time was 2:18 pm" (should have been 2:04 pm, so the recording was started 34 minutes after treatment was added rather than the usual 20 minutes)'
$ws.Range("J7").Value = 'Regarding Plate 72-8208, Seline explained that after hitting start, Axion runs an “offset” as a preliminary calibration step. I guess they usually hit record after this offset has completed. However, the offset only takes about 10 seconds. Therefore, for a 40 minute recording, the effect of including a 10-sec offset should probably not affect the parameter values too much.'
$ws.Range("K7").Value = 'Will see if plate 72-8209 looks that much different than the other replicates (given the decline in activity over time, an extra 14 mintues of ''rest'' after treatment addition might matter?)'

# Row 20 (culture G19): clarify that well quality was already zeroed for G3/H3.
$ws.Range("K20").Value = 'check if plate 75-8114 looks okay despite tipping. Well quality already set to 0 for wells G3 and H3, but check if other wells appear usable'

# Rows 23, 25 and 33: reword the LDH/AB sinking follow-up note.
$ws.Range("J23").Value = 'See convo with Tim 3/7/23 - MEA activity should not be affected, but there could be interference in the LDH or AB. Do follow up analysis). Still need to check up on whether we care about solution sinking to bottom of well. '
$ws.Range("J25").Value = 'See convo with Tim 3/7/23 - MEA activity should not be affected, but there could be interference in the LDH or AB. Do follow up analysis). Still need to check up on whether we care about solution sinking to bottom of well, still need to do follow up on combining recording (though that might be analysis)'
$ws.Range("J33").Value = 'See convo with Tim 3/7/23 - MEA activity should not be affected, but there could be interference in the LDH or AB. Do follow up analysis). Still need to check up on whether we care about solution sinking to bottom of well, still need to do follow up on combining recording (though that might be analysis)'

# Update the sheet view: scroll the frozen pane back to row 2 and move the
# active selection to K7.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
[void]$ws.Range("K7").Select()
